# Swap the contents of column C (codeforiati:group-code) and column D
# (codeforiati:group-name), including the header row, for every used row
# on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row   # xlUp = -4162

$rangeC = $ws.Range("C1:C$lastRow")
$rangeD = $ws.Range("D1:D$lastRow")

$valuesC = $rangeC.Value2
$valuesD = $rangeD.Value2

$rangeC.Value2 = $valuesD
$rangeD.Value2 = $valuesC
